$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44574
$ws.Range("K2").Value = 'Black Amber'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18500
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'
$ws.Range("S2").Value = 1028

# Row 3
$ws.Range("D3").Value = 44650
$ws.Range("K3").Value = 'Angeleno'
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17500
$ws.Range("Q3").Value = '$/bandeja 18 kilos granel'
$ws.Range("S3").Value = 972

# Row 4
$ws.Range("D4").Value = 44596
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("S4").Value = 861

# Row 5
$ws.Range("D5").Value = 44169
$ws.Range("K5").Value = 'Angeleno'
$ws.Range("L5").Value = 'Tercera'
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("S5").Value = 1361

# Row 7
$ws.Range("D7").Value = 44174
$ws.Range("K7").Value = 'Angeleno'
$ws.Range("L7").Value = 'Primera'
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20500
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("S7").Value = 1139

# Row 8
$ws.Range("D8").Value = 44175
$ws.Range("K8").Value = 'Angeleno'
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 21000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21500
$ws.Range("S8").Value = 1194

# Row 9
$ws.Range("D9").Value = 44706
$ws.Range("K9").Value = 'Angeleno'
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15500
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 861

# Row 10
$ws.Range("D10").Value = 44314
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("Q10").Value = '$/bandeja 18 kilos granel'
$ws.Range("S10").Value = 806

# Row 11
$ws.Range("D11").Value = 44243
$ws.Range("K11").Value = 'Black Amber'
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 14500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("S11").Value = 806

# Row 12
$ws.Range("D12").Value = 44587
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15500
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("S12").Value = 861

# Row 13
$ws.Range("D13").Value = 44239
$ws.Range("L13").Value = 'Primera'
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15500
$ws.Range("S13").Value = 861

# Row 14
$ws.Range("D14").Value = 44217
$ws.Range("K14").Value = 'Black Amber'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 17000
$ws.Range("P14").Value = 16500
$ws.Range("S14").Value = 917

# Row 15
$ws.Range("D15").Value = 44580
$ws.Range("K15").Value = 'Black Amber'
$ws.Range("M15").Value = 270
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19500
$ws.Range("R15").Value = 'Región Metropolitana'
$ws.Range("S15").Value = 1083

# Row 16
$ws.Range("D16").Value = 44614
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 19000
$ws.Range("P16").Value = 18500
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 1028

# Row 17
$ws.Range("D17").Value = 44278
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 15500
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("S17").Value = 861

# Row 18
$ws.Range("D18").Value = 44245
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 250
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 14500
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 806

# Row 19
$ws.Range("D19").Value = 44229
$ws.Range("K19").Value = 'Fortuna'
$ws.Range("N19").Value = 14000
$ws.Range("O19").Value = 15000
$ws.Range("P19").Value = 14500
$ws.Range("S19").Value = 806

# Row 20
$ws.Range("D20").Value = 44628
$ws.Range("M20").Value = 270
$ws.Range("Q20").Value = '$/bandeja 18 kilos granel'

# Row 21
$ws.Range("D21").Value = 44238
$ws.Range("K21").Value = 'Black Amber'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("N21").Value = 14000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 14500
$ws.Range("S21").Value = 806

# Row 22
$ws.Range("D22").Value = 44238
